$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.788.50'
$ws.Range("E2").Value = '  +2.98%  '

$ws.Range("D3").Value = '1.678.09'
$ws.Range("E3").Value = '  +2.97%  '

$ws.Range("D4").Formula = "'0.998"
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").Formula = "'219.13"
$ws.Range("E5").Value = '  +2.21%  '

$ws.Range("D6").Formula = "'0.529"
$ws.Range("E6").Value = '  +2.18%  '

$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("D8").Formula = "'29.17"
$ws.Range("E8").Value = '  +2.48%  '

$ws.Range("E9").Value = '  +2.18%  '

$ws.Range("D10").Formula = "'0.0645"
$ws.Range("E10").Value = '  +6.13%  '

$ws.Range("D11").Formula = "'0.0903"
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("D12").Value = '1.918.63'
$ws.Range("E12").Value = '  +2.88%  '

$ws.Range("D13").Value = '1.675.08'
$ws.Range("E13").Value = '  +2.60%  '

$ws.Range("D14").Formula = "'10.15"
$ws.Range("E14").Value = '  +9.80%  '

$ws.Range("E15").Value = '  +7.96%  '

$ws.Range("D16").Formula = "'4.03"
$ws.Range("E16").Value = '  +5.25%  '

$ws.Range("D17").Value = '30.751.65'
$ws.Range("E17").Value = '  +2.81%  '

$ws.Range("D18").Formula = "'65.93"
$ws.Range("E18").Value = '  +3.23%  '

$ws.Range("D19").Formula = "'243.67"
$ws.Range("E19").Value = '  +1.31%  '

$ws.Range("D20").Value = '0.0₃0721'
$ws.Range("E20").Value = '  +3.22%  '

$ws.Range("D21").Formula = "'0.998"
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("E22").Value = '  +2.80%  '

$ws.Range("D23").Formula = "'9.98"
$ws.Range("E23").Value = '  +2.19%  '

$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").Formula = "'159.31"
$ws.Range("E25").Value = '  +0.98%  '

$ws.Range("D26").Formula = "'15.80"
$ws.Range("E26").Value = '  +2.46%  '

$ws.Range("E27").Value = '  +2.49%  '

$ws.Range("E28").Value = '  +2.19%  '

$ws.Range("D29").Formula = "'0.999"
$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("E30").Value = '  +1.35%  '

$ws.Range("E31").Value = '  +3.86%  '

$ws.Range("D32").Formula = "'3.46"
$ws.Range("E32").Value = '  +2.91%  '

$ws.Range("D33").Value = '1.519.17'
$ws.Range("E33").Value = '  +6.68%  '

$ws.Range("D34").Formula = "'3.31"
$ws.Range("E34").Value = '  +4.57%  '

$ws.Range("E35").Value = '  +6.26%  '

$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("D37").Formula = "'83.27"
$ws.Range("E37").Value = '  +11.43%  '

$ws.Range("D38").Formula = "'0.601"
$ws.Range("E38").Value = '  +8.52%  '

$ws.Range("E39").Value = '  +4.74%  '

$ws.Range("D40").Formula = "'2.65"
$ws.Range("E40").Value = '  -3.21%  '

$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").Formula = "'2.02"
$ws.Range("E42").Value = '  +2.20%  '

$ws.Range("D43").Formula = "'0.837"
$ws.Range("E43").Value = '  +1.42%  '

$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("E45").Value = '  +2.01%  '

$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("D47").Formula = "'5.53"
$ws.Range("E47").Value = '  +3.81%  '

$ws.Range("D48").Value = '1.811.70'
$ws.Range("E48").Value = '  +2.22%  '

$ws.Range("D49").Formula = "'50.40"
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0116'
$ws.Range("E50").Value = '  +4.38%  '

$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Formula = "'92.60"
$ws.Range("E51").Value = '  +2.16%  '
